$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61, shifting existing rows 61-122 down to 62-123
$ws.Rows.Item(61).Insert()

# Fill in the new row 61 with data
$ws.Cells.Item(61, 1).Value = 4
$ws.Cells.Item(61, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(61, 3).Value = "Los Lagos"
$ws.Cells.Item(61, 4).Value = 44484
$ws.Cells.Item(61, 5).Value = 10
$ws.Cells.Item(61, 6).Value = 100112039
$ws.Cells.Item(61, 7).Value = "Ciboulette"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 240
$ws.Cells.Item(61, 11).Value = 2500
$ws.Cells.Item(61, 12).Value = 2500
$ws.Cells.Item(61, 13).Value = 2500
$ws.Cells.Item(61, 14).Value = "$/docena de atados"
$ws.Cells.Item(61, 15).Value = "Región Metropolitana"
$ws.Cells.Item(61, 16).Value = 833
$ws.Cells.Item(61, 17).Value = 3
$ws.Cells.Item(61, 18).Value = "Hortaliza"
